$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taxon")

# Insert a new column before column G (taxonomicNodeId's right neighbour /
# title's current position) to make room, shifting title..iri one column right.
$ws.Columns("G:G").Insert()

# previouslyKnownAs used to live in D1; it is being replaced there by the
# column that used to be E1 (externalEquivalentTaxon). Shift D1:F1 left by one
# (drop the old previouslyKnownAs), then place the relocated previouslyKnownAs
# and the new alternateName column in G1:H1.
$ws.Range("D1").Value = "externalEquivalentTaxon"
$ws.Range("E1").Value = "taxonomicId"
$ws.Range("F1").Value = "taxonomicNodeId"
$ws.Range("G1").Value = "alternateName"
$ws.Range("H1").Value = "previouslyKnownAs"
